$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force columns D (Price) and E (Volume) to remain plain text so that
# numeric-looking strings (e.g. "0.9985", "29.019.61") are not silently
# re-interpreted by Excel as numbers/dates, matching the source data which
# is entirely textual (prices with thousands separators as dots, etc).
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '29.019.61'
$ws.Range("E2").Value = '  +5.74%  '

$ws.Range("D3").Value = '1.916.52'
$ws.Range("E3").Value = '  +5.10%  '

$ws.Range("D4").Value = '0.9985'

$ws.Range("D5").Value = '338.74'
$ws.Range("E5").Value = '  +1.90%  '

$ws.Range("D6").Value = '0.9988'
$ws.Range("E6").Value = '  -0.58%  '

$ws.Range("D7").Value = '0.4742'
$ws.Range("E7").Value = '  +3.78%  '

$ws.Range("D8").Value = '0.4063'
$ws.Range("E8").Value = '  +7.08%  '

$ws.Range("D9").Value = '48.01'
$ws.Range("E9").Value = '  +3.86%  '

$ws.Range("D10").Value = '0.08196'
$ws.Range("E10").Value = '  +4.01%  '

$ws.Range("D11").Value = '1.035'
$ws.Range("E11").Value = '  +6.95%  '

$ws.Range("D12").Value = '22.58'
$ws.Range("E12").Value = '  +7.64%  '

$ws.Range("B13").Value = 'WrappedEther'
$ws.Range("C13").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D13").Value = '1.903.32'
$ws.Range("E13").Value = '  +4.87%  '

$ws.Range("B14").Value = 'Polkadot'
$ws.Range("C14").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D14").Value = '6.114'
$ws.Range("E14").Value = '  +4.09%  '

$ws.Range("D15").Value = '7.408'
$ws.Range("E15").Value = '  +5.30%  '

$ws.Range("D16").Value = '91.70'
$ws.Range("E16").Value = '  +2.35%  '

$ws.Range("D17").Value = '0.9993'
$ws.Range("E17").Value = '  -0.69%  '

$ws.Range("D18").Value = '0.00001055'
$ws.Range("E18").Value = '  +3.00%  '

$ws.Range("D19").Value = '0.06621'
$ws.Range("E19").Value = '  -0.28%  '

$ws.Range("D20").Value = '17.97'
$ws.Range("E20").Value = '  +5.35%  '

$ws.Range("D21").Value = '0.9993'
$ws.Range("E21").Value = '  -0.49%  '

$ws.Range("D22").Value = '29.036.41'
$ws.Range("E22").Value = '  +5.89%  '

$ws.Range("D23").Value = '5.579'
$ws.Range("E23").Value = '  +4.79%  '

$ws.Range("D24").Value = '11.23'
$ws.Range("E24").Value = '  +4.11%  '

$ws.Range("D25").Value = '2.268'
$ws.Range("E25").Value = '  -1.73%  '

$ws.Range("D26").Value = '2.130.12'
$ws.Range("E26").Value = '  +4.87%  '

$ws.Range("D27").Value = '160.83'
$ws.Range("E27").Value = '  +3.44%  '

$ws.Range("D28").Value = '20.06'
$ws.Range("E28").Value = '  +3.77%  '

$ws.Range("D29").Value = '2.190'
$ws.Range("E29").Value = '  +6.91%  '

$ws.Range("D30").Value = '5.568'
$ws.Range("E30").Value = '  +5.67%  '

$ws.Range("D31").Value = '121.19'
$ws.Range("E31").Value = '  +2.63%  '

$ws.Range("D32").Value = '1.023'
$ws.Range("E32").Value = '  +8.96%  '

$ws.Range("D33").Value = '0.09592'
$ws.Range("E33").Value = '  +3.24%  '

$ws.Range("B34").Value = 'ARBITRUM'
$ws.Range("C34").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D34").Value = '1.436'
$ws.Range("E34").Value = '  +9.23%  '

$ws.Range("B35").Value = 'HuobiToken'
$ws.Range("C35").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D35").Value = '3.655'
$ws.Range("E35").Value = '  +2.11%  '

$ws.Range("D36").Value = '5.442'
$ws.Range("E36").Value = '  +3.91%  '

$ws.Range("D37").Value = '0.06229'
$ws.Range("E37").Value = '  +5.16%  '

$ws.Range("D38").Value = '0.02292'
$ws.Range("E38").Value = '  +5.28%  '

$ws.Range("D39").Value = '8.712'
$ws.Range("E39").Value = '  +8.41%  '

$ws.Range("D40").Value = '1.208'
$ws.Range("E40").Value = '  +5.44%  '

$ws.Range("D41").Value = '0.6068'
$ws.Range("E41").Value = '  +5.49%  '

$ws.Range("B42").Value = 'Aptos'
$ws.Range("C42").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D42").Value = '10.65'
$ws.Range("E42").Value = '  +7.14%  '

$ws.Range("B43").Value = 'Algorand'
$ws.Range("C43").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D43").Value = '0.1909'
$ws.Range("E43").Value = '  +4.80%  '

$ws.Range("D44").Value = '0.9983'
$ws.Range("E44").Value = '  -0.53%  '

$ws.Range("E45").Value = '  +0.15%  '

$ws.Range("B46").Value = 'Decentraland'
$ws.Range("C46").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D46").Value = '0.5662'
$ws.Range("E46").Value = '  +4.17%  '

$ws.Range("B47").Value = 'EnergySwap'
$ws.Range("C47").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D47").Value = '12.38'
$ws.Range("E47").Value = '  +3.48%  '

$ws.Range("E48").Value = '  +6.91%  '

$ws.Range("D49").Value = '0.07307'
$ws.Range("E49").Value = '  +10.77%  '

$ws.Range("D50").Value = '2.169'
$ws.Range("E50").Value = '  +19.74%  '

$ws.Range("D51").Value = '113.35'
$ws.Range("E51").Value = '  +2.53%  '
